{"js": "// Find the paragraph that contains the sentence about \"agentName\" and\n// make the word \"agentName\" bold, splitting the single run into three\n// runs: \"The \", \"agentName\" (bold), \" is a general text string giving\n// the agent a name.\"\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst target = paragraphs.items.find(\n  (p) => p.text === \"The agentName is a general text string giving the agent a name.\"\n);\n\nif (!target) {\n  throw new Error(\"Could not find the target paragraph.\");\n}\n\nconst results = target.search(\"agentName\", { matchCase: true, matchWholeWord: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find 'agentName' inside the target paragraph.\");\n}\n\nresults.items[0].font.bold = true;\nawait context.sync();\n", "ps1": "# Find the paragraph that contains the sentence about \"agentName\" and\n# make the word \"agentName\" bold, splitting the single run into three\n# runs: \"The \", \"agentName\" (bold), \" is a general text string giving\n# the agent a name.\"\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text\n    if ($text -like \"The agentName is a general text string giving the agent a name.*\") {\n        $r = $p.Range\n        $r.Find.ClearFormatting()\n        $r.Find.MatchCase = $true\n        $r.Find.MatchWholeWord = $true\n        $found = $r.Find.Execute(\"agentName\")\n        if ($found) {\n            $r.Font.Bold = 1\n        }\n        break\n    }\n}\n"}
